# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Leve profit tracker sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 27367.578
$ws.Range("J21").Value = 26110.223
$ws.Range("L21").Value = 26110.223
$ws.Range("N21").Value = -27046.223

$ws.Range("H23").Value = 27367.578
$ws.Range("J23").Value = 26110.223
$ws.Range("L23").Value = 26110.223
$ws.Range("N23").Value = -26578.223

$ws.Range("H40").Value = 2283.3333
$ws.Range("I40").Value = 2250
$ws.Range("K40").Value = 2250
$ws.Range("M40").Value = -2075

$ws.Range("H62").Value = 42453.88
$ws.Range("I62").Value = 48807
$ws.Range("K62").Value = 48807
$ws.Range("M62").Value = -48183

$ws.Range("H65").Value = 42453.88
$ws.Range("I65").Value = 48807
$ws.Range("K65").Value = 244035
$ws.Range("M65").Value = -240915

$ws.Range("H135").Value = 9068.143
$ws.Range("I135").Value = 9068.143
$ws.Range("K135").Value = 81613.287
$ws.Range("M135").Value = -79078.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14855.156
$ws.Range("I32").Value = 1524.381
$ws.Range("K32").Value = 1524.381
$ws.Range("M32").Value = -1237.381

$ws.Range("H43").Value = 8450
$ws.Range("J43").Value = 8450
$ws.Range("L43").Value = 8450
$ws.Range("N43").Value = -9076

$ws.Range("H133").Value = 32083.334
$ws.Range("J133").Value = 32083.334
$ws.Range("L133").Value = 32083.334
$ws.Range("N133").Value = -37143.334

$ws.Range("H139").Value = 38750
$ws.Range("J139").Value = 38750
$ws.Range("L139").Value = 38750
$ws.Range("N139").Value = -49030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 12398
$ws.Range("J132").Value = 12398
$ws.Range("L132").Value = 12398
$ws.Range("N132").Value = -22518

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1344.762
$ws.Range("I31").Value = 1340.2941
$ws.Range("J31").Value = 1363.75
$ws.Range("K31").Value = 1340.2941
$ws.Range("L31").Value = 1363.75
$ws.Range("M31").Value = -1045.2941
$ws.Range("N31").Value = -1953.75

$ws.Range("H34").Value = 1344.762
$ws.Range("I34").Value = 1340.2941
$ws.Range("J34").Value = 1363.75
$ws.Range("K34").Value = 1340.2941
$ws.Range("L34").Value = 1363.75
$ws.Range("M34").Value = -1138.2941
$ws.Range("N34").Value = -1767.75

$ws.Range("H59").Value = 32450
$ws.Range("J59").Value = 32450
$ws.Range("L59").Value = 32450
$ws.Range("N59").Value = -34740

$ws.Range("H132").Value = 2805.2
$ws.Range("I132").Value = 2139.4285
$ws.Range("J132").Value = 4358.6665
$ws.Range("K132").Value = 6418.2855
$ws.Range("L132").Value = 13075.9995
$ws.Range("M132").Value = -3888.2855
$ws.Range("N132").Value = -18135.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 751
$ws.Range("J17").Value = 751
$ws.Range("L17").Value = 2253
$ws.Range("N17").Value = -2591

$ws.Range("H22").Value = 2237.5
$ws.Range("J22").Value = 2485.7144
$ws.Range("L22").Value = 7457.1432
$ws.Range("N22").Value = -7795.1432

$ws.Range("H27").Value = 2237.5
$ws.Range("J27").Value = 2485.7144
$ws.Range("L27").Value = 7457.1432
$ws.Range("N27").Value = -7661.1432

$ws.Range("H46").Value = 628.6
$ws.Range("I46").Value = 47.666668
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 143.000004
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -52.00000399999999
$ws.Range("N46").Value = -4682

$ws.Range("H63").Value = 3533.3333
$ws.Range("I63").Value = 3800
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 11400
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -10651
$ws.Range("N63").Value = -10498

$ws.Range("H66").Value = 3533.3333
$ws.Range("I66").Value = 3800
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 34200
$ws.Range("L66").Value = 27000
$ws.Range("M66").Value = -30456
$ws.Range("N66").Value = -34488

$ws.Range("H69").Value = 4499.727
$ws.Range("J69").Value = 4499.727
$ws.Range("L69").Value = 13499.181
$ws.Range("N69").Value = -15121.181

$ws.Range("H72").Value = 4499.727
$ws.Range("J72").Value = 4499.727
$ws.Range("L72").Value = 40497.543
$ws.Range("N72").Value = -48609.543

$ws.Range("H129").Value = 2280.5
$ws.Range("I129").Value = 2360
$ws.Range("J129").Value = 2229.9092
$ws.Range("K129").Value = 7080
$ws.Range("L129").Value = 6689.7276
$ws.Range("M129").Value = -2080
$ws.Range("N129").Value = -16689.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8975
$ws.Range("I80").Value = 14875
$ws.Range("J80").Value = 3075
$ws.Range("K80").Value = 14875
$ws.Range("L80").Value = 3075
$ws.Range("M80").Value = -13877
$ws.Range("N80").Value = -5071

$ws.Range("H83").Value = 8975
$ws.Range("I83").Value = 14875
$ws.Range("J83").Value = 3075
$ws.Range("K83").Value = 74375
$ws.Range("L83").Value = 15375
$ws.Range("M83").Value = -69383
$ws.Range("N83").Value = -25359

$ws.Range("H103").Value = 20300
$ws.Range("J103").Value = 20300
$ws.Range("L103").Value = 20300
$ws.Range("N103").Value = -22644

$ws.Range("H138").Value = 52999.668
$ws.Range("J138").Value = 52999.668
$ws.Range("L138").Value = 52999.668
$ws.Range("N138").Value = -63279.668

$ws.Range("H139").Value = 46000
$ws.Range("J139").Value = 46000
$ws.Range("L139").Value = 46000
$ws.Range("N139").Value = -56280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3348.3333
$ws.Range("I7").Value = 2973.5715
$ws.Range("J7").Value = 3535.7144
$ws.Range("K7").Value = 2973.5715
$ws.Range("L7").Value = 3535.7144
$ws.Range("M7").Value = -2861.5715
$ws.Range("N7").Value = -3759.7144

$ws.Range("H22").Value = 13489
$ws.Range("I22").Value = 593.3333
$ws.Range("J22").Value = 21226.4
$ws.Range("K22").Value = 593.3333
$ws.Range("L22").Value = 21226.4
$ws.Range("M22").Value = -298.3333
$ws.Range("N22").Value = -21816.4

$ws.Range("H27").Value = 13489
$ws.Range("I27").Value = 593.3333
$ws.Range("J27").Value = 21226.4
$ws.Range("K27").Value = 593.3333
$ws.Range("L27").Value = 21226.4
$ws.Range("M27").Value = -486.3333
$ws.Range("N27").Value = -21440.4

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null

$ws.Range("H93").Value = 519.3333
$ws.Range("I93").Value = 496.45834
$ws.Range("J93").Value = 610.8333
$ws.Range("K93").Value = 496.45834
$ws.Range("L93").Value = 610.8333
$ws.Range("M93").Value = 751.54166
$ws.Range("N93").Value = -3106.8333

$ws.Range("H126").Value = 3348.3333
$ws.Range("I126").Value = 2973.5715
$ws.Range("J126").Value = 3535.7144
$ws.Range("K126").Value = 8920.7145
$ws.Range("L126").Value = 10607.1432
$ws.Range("M126").Value = -6450.7145
$ws.Range("N126").Value = -15547.1432

$ws.Range("H132").Value = 19999.75
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
